$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-21 06:48:41"
$wsZh.Range("H4").Value = "2016-03-21 06:49:26"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-21 06:48:49"
$wsDe.Range("H4").Value = "2016-03-21 06:49:43"
